$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 19 data: a new leetcode entry about removing zero-sum consecutive
# nodes from a linked list. Shared-string insertion order matters (matches the
# authoring order: method, question, keywords, then space complexity), so set
# the cells in that sequence.
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = 1171
$ws.Range("D19").Value = "solder方便处理1-2--3`n1 使用map，key是当前节点以及之前节点的累加和，value是当前节点对象`n2 迭代节点，cur对应的节点和存在多次，说明这中间是有连续元素和为0的节点，即1-2-3-4,<0,0><1,1><3,2><6,3><3,-3><7,4>`n3 找到这个节点cur1，将cur.next=cur1.next`n4 如果没有连续节点和是0，在map中也能找到一次累加和，执行cur.next=cur1.next，也是一样的，因为cur=cur1"
$ws.Range("C19").Value = "给你一个链表的头节点 head，请你编写代码，反复删去链表中由 总和 值为 0 的连续节点组成的序列，直到不存在这样的序列为止"
$ws.Range("E19").Value = "solder`n删除节点`n节点累加和`n哈希表"
$ws.Range("F19").Value = "O(N), N是元素个数"
$ws.Range("G19").Value = "O(K)`nK是链表节点个数"

# Match row height used by similarly-sized rows (row 15 uses 200 too).
$ws.Rows.Item(19).RowHeight = 200

# Update the sheet's active selection, matching the new content layout.
$ws.Range("C22").Select()
